# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# This updates the "K" column (column G) values for each data row on the
# active sheet to the newly-computed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K value (column G)
$kValues = @{
    2  = 4
    3  = 1
    4  = 1
    5  = 1
    6  = 2
    7  = 1
    8  = 1
    9  = 0
    10 = 1
    11 = 3
    12 = 2
    13 = 1
    14 = 2
    15 = 1
    16 = 0
    17 = 1
    18 = 0
    19 = 1
    20 = 0
    21 = 0
    22 = 1
    23 = 1
    24 = 3
    25 = 0
    26 = 2
    27 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
